$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Wipe all existing hyperlinks (and their relationships) so we can re-add
#    them in a clean, deterministic rId1..rId14 order that matches the new
#    row layout.
# ---------------------------------------------------------------------------
$ws.Range("B2").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2) New header cell F1 = "REX_DEF" (bold/boxed header style, like B1:E1)
# ---------------------------------------------------------------------------
$ws.Range("F1").Value2 = "REX_DEF"
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Row data (A=index, B=OBI IRI, C=OBI desc dict, D=REX IRI, E=REX desc
#    dict, F=REX def) for the 7 data rows (rows 2-8).
# ---------------------------------------------------------------------------
$rows = @(
  @{ Row=2; A=0; B="http://purl.obolibrary.org/obo/OBI_0000374"; C="{'label': 'excitation', 'prefLabel': None, 'altLabel': None, 'name': 'OBI_0000374'}"; D="http://purl.obolibrary.org/obo/REX_0000026"; E="{'label': 'excitation'}"; F="[]" },
  @{ Row=3; A=1; B="http://purl.obolibrary.org/obo/OBI_0600038"; C="{'label': 'addition', 'prefLabel': None, 'altLabel': None, 'name': 'OBI_0600038'}"; D="http://purl.obolibrary.org/obo/REX_0000089"; E="{'label': 'addition'}"; F="[]" },
  @{ Row=4; A=2; B="http://purl.obolibrary.org/obo/OBI_0000213"; C="{'label': 'fluorescence', 'prefLabel': None, 'altLabel': None, 'name': 'OBI_0000213'}"; D="http://purl.obolibrary.org/obo/REX_0000043"; E="{'label': 'fluorescence'}"; F="[]" },
  @{ Row=5; A=3; B="http://purl.obolibrary.org/obo/OBI_0600034"; C="{'label': 'precipitation', 'prefLabel': None, 'altLabel': None, 'name': 'OBI_0600034'}"; D="http://purl.obolibrary.org/obo/REX_0000182"; E="{'label': 'precipitation'}"; F="[]" },
  @{ Row=6; A=4; B="http://purl.obolibrary.org/obo/OBI_0000385"; C="{'label': 'ionization', 'prefLabel': None, 'altLabel': None, 'name': 'OBI_0000385'}"; D="http://purl.obolibrary.org/obo/REX_0000152"; E="{'label': 'ionization'}"; F="[]" },
  @{ Row=7; A=5; B="http://purl.obolibrary.org/obo/OBI_0302890"; C="{'label': 'polymerization', 'prefLabel': None, 'altLabel': None, 'name': 'OBI_0302890'}"; D="http://purl.obolibrary.org/obo/REX_0000251"; E="{'label': 'polymerization'}"; F="[]" },
  @{ Row=8; A=6; B="http://purl.obolibrary.org/obo/OBI_0600053"; C="{'label': 'electrophoresis', 'prefLabel': None, 'altLabel': None, 'name': 'OBI_0600053'}"; D="http://purl.obolibrary.org/obo/REX_0000338"; E="{'label': 'electrophoresis'}"; F="[]" }
)

# ---------------------------------------------------------------------------
# 4) Write the cell values for every data row.
# ---------------------------------------------------------------------------
foreach ($r in $rows) {
  $n = $r.Row
  $ws.Range("A$n").Value2 = $r.A
  $ws.Range("B$n").Value2 = $r.B
  $ws.Range("C$n").Value2 = $r.C
  $ws.Range("D$n").Value2 = $r.D
  $ws.Range("E$n").Value2 = $r.E
  $ws.Range("F$n").Value2 = $r.F
}

# ---------------------------------------------------------------------------
# 5) Re-create hyperlinks in B/D columns, in row order, so relationship ids
#    come out as rId1..rId14 (matching the target). Reset the B:D area to
#    the plain/default style first so every cell starts from the same
#    baseline before Hyperlinks.Add applies its own auto-format - this keeps
#    the number of distinct (and ultimately unused, once we repaint in step
#    6) style variants minted by the auto-formatting down to a minimum.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("B2:D8").PasteSpecial(-4122)

foreach ($r in $rows) {
  $n = $r.Row
  $ws.Hyperlinks.Add($ws.Range("B$n"), $r.B)
  $ws.Hyperlinks.Add($ws.Range("D$n"), $r.D)
}

# ---------------------------------------------------------------------------
# 6) Re-apply the correct cell formatting (Hyperlinks.Add overwrites it with
#    its own auto-format, so fix up after adding links):
#      - column A: same style as A2 (bold, boxed, centered)
#      - columns B & D: same style as B2 (Hyperlink style)
#      - columns C, E, F: default/no style
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B2:B8").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D2:D8").PasteSpecial(-4122)

Write-Output "Mappings updated"
